$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# 1) Duplicate "2022-Q1" sheet (placed right after original) so we keep a style-1-correct copy.
$ws2.Copy($null, $ws2)
$ws3 = $wb.Worksheets.Item(3)

# 2) Rename: original becomes "2022-Q3" (keeps sheetId/rId of old "2022-Q1"); copy becomes "2022-Q1".
$ws2.Name = "2022-Q3"
$ws3.Name = "2022-Q1"

# 3) Clear old data out of the sheet that is now "2022-Q3": wipe rows 1-4 in place, and
#    altogether remove the now-unneeded rows 5-8 so the used range shrinks back down.
$ws2.Range("A1:H4").Clear()
$ws2.Range("A5:H8").EntireRow.Delete()

# 4) Re-apply the "style 2" look (bold/border/center, same as the 总计 header row) used by the
#    Q3 sheet's header row + index column, by copying format from sheet1 which already has it.
$ws1.Range("B1:D1").Copy()
$ws2.Range("B1:H1").PasteSpecial(-4122)
$ws1.Range("A2").Copy()
$ws2.Range("A2:A4").PasteSpecial(-4122)

# 5) Header row text.
$ws2.Range("B1").Value = "基金代码"
$ws2.Range("C1").Value = "基金名称"
$ws2.Range("D1").Value = "基金规模"
$ws2.Range("E1").Value = "股票总仓位"
$ws2.Range("F1").Value = "仓位占比"
$ws2.Range("G1").Value = "持有市值(亿元)"
$ws2.Range("H1").Value = "仓位排名"

# 6) Data rows. B and D:G must stay text (numeric-looking strings): force text format so the
#    assignment isn't silently coerced to a number, then strip the format back off again (the
#    target file leaves these cells with no explicit style -- only the stored type matters).
$ws2.Range("B2:B4").NumberFormat = "@"
$ws2.Range("D2:G4").NumberFormat = "@"

$ws2.Range("A2").Value = 0
$ws2.Range("B2").Value = "160143"
$ws2.Range("C2").Value = "南方创业板2年定期开放混合"
$ws2.Range("D2").Value = "3.38"
$ws2.Range("E2").Value = "83.97"
$ws2.Range("F2").Value = "4.17"
$ws2.Range("G2").Value = "0.1409"
$ws2.Range("H2").Value = 2

$ws2.Range("A3").Value = 1
$ws2.Range("B3").Value = "002160"
$ws2.Range("C3").Value = "南方转型驱动灵活配置混合"
$ws2.Range("D3").Value = "3.14"
$ws2.Range("E3").Value = "93.43"
$ws2.Range("F3").Value = "2.40"
$ws2.Range("G3").Value = "0.0754"
$ws2.Range("H3").Value = 7

$ws2.Range("A4").Value = 2
$ws2.Range("B4").Value = "003513"
$ws2.Range("C4").Value = "中邮消费升级灵活配置混合"
$ws2.Range("D4").Value = "0.56"
$ws2.Range("E4").Value = "30.56"
$ws2.Range("F4").Value = "3.04"
$ws2.Range("G4").Value = "0.0170"
$ws2.Range("H4").Value = 5

$ws2.Range("B2:B4").ClearFormats()
$ws2.Range("D2:G4").ClearFormats()

# 7) Update the "总计" summary sheet: insert a new row 2 for 2022-Q3, push 2022-Q1's row to row 3.
$ws1.Range("A2:D2").Copy()
$ws1.Range("A3:D3").PasteSpecial(-4122)
$ws1.Range("A3").Value = 1
$ws1.Range("B3").Value = "2022-Q1"
$ws1.Range("C3").Value = 7
$ws1.Range("D3").Value = 0.38

$ws1.Range("A2").Value = 0
$ws1.Range("B2").Value = "2022-Q3"
$ws1.Range("C2").Value = 3
$ws1.Range("D2").Value = 0.23

Write-Output "done"
